# Update CDRDfRCP parameters to not start increasing until ~95% CES
$wb = $excel.ActiveWorkbook

$wsCDR = $wb.Worksheets.Item("CDRDfRCP")
$wsAbout = $wb.Worksheets.Item("About")

# --- Core parameter edits (a parameter 15 -> 75, b parameter 0.9 -> 0.98) ---
# These drive the (1-EXP(-((A)/B2)^B1))*B3 curve on the About sheet, which
# recalculates automatically, and the chart's cached values along with it.
$wsCDR.Range("B1").Value = 75
$wsCDR.Range("B2").Value = 0.98

# --- Cosmetic: widen column K on the About sheet (used for the reviewer's
# new calibration notes off to the side of the data table) ---
$wsAbout.Columns.Item(11).ColumnWidth = 19.3

# --- Move/resize the chart on the About sheet to its new anchor position ---
$co = $wsAbout.ChartObjects(1)
$co.Left = 201.5625
$co.Top = 172.87496062992125
$co.Width = 443.5
$co.Height = 216

# --- Selections left behind by the editing session ---
# CDRDfRCP is visited and left with C29 selected, but About remains the
# active (front-most) tab, so activate it last.
$wsCDR.Activate()
$wsCDR.Range("C29").Select()

$wsAbout.Activate()
$wsAbout.Range("K2:L4").Select()
